# Apply the dated worksheet updates: refresh the date header and all
# "three-digit number divided by one-digit number" problems.

$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "2024-02-18 Sunday"; New = "2024-02-19 Monday" },
    @{ Old = "571÷5="; New = "295÷4=" },
    @{ Old = "856÷9="; New = "234÷7=" },
    @{ Old = "218÷4="; New = "939÷9=" },
    @{ Old = "922÷3="; New = "204÷3=" },
    @{ Old = "935÷5="; New = "587÷4=" },
    @{ Old = "609÷2="; New = "291÷4=" },
    @{ Old = "174÷5="; New = "517÷9=" },
    @{ Old = "784÷2="; New = "317÷4=" },
    @{ Old = "176÷7="; New = "130÷5=" },
    @{ Old = "175÷3="; New = "745÷8=" },
    @{ Old = "151÷7="; New = "454÷5=" },
    @{ Old = "548÷7="; New = "406÷3=" },
    @{ Old = "663÷9="; New = "163÷8=" },
    @{ Old = "143÷4="; New = "851÷5=" },
    @{ Old = "971÷3="; New = "894÷9=" },
    @{ Old = "177÷6="; New = "445÷4=" },
    @{ Old = "662÷3="; New = "501÷3=" },
    @{ Old = "103÷7="; New = "473÷7=" },
    @{ Old = "252÷3="; New = "868÷4=" },
    @{ Old = "188÷6="; New = "623÷8=" },
    @{ Old = "953÷5="; New = "172÷4=" },
    @{ Old = "207÷7="; New = "482÷2=" },
    @{ Old = "844÷9="; New = "963÷6=" },
    @{ Old = "601÷8="; New = "896÷8=" },
    @{ Old = "680÷4="; New = "855÷4=" }
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.Old, $true, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)
}
